$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header A1: "value" -> "Protection"
$ws.Range("A1").Value = "Protection"

# A2: remove internal line breaks, join with single spaces
$ws.Range("A2").Value = "Insufficient legal protection"

# A4: remove internal line break, join with single space
$ws.Range("A4").Value = "Production areas"

# A6: remove internal line breaks, join with single spaces
$ws.Range("A6").Value = "Requires individual assessment"
